$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171, shifting the existing rows 171..205 down to 172..206
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new weekly price record
$ws.Cells.Item(171, 1).Value = 3
$ws.Cells.Item(171, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44476
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 6).Value = 100112009
$ws.Cells.Item(171, 7).Value = "Acelga"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 270
$ws.Cells.Item(171, 11).Value = 2000
$ws.Cells.Item(171, 12).Value = 2200
$ws.Cells.Item(171, 13).Value = 2089
$ws.Cells.Item(171, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(171, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(171, 16).Value = 348
$ws.Cells.Item(171, 17).Value = 6
$ws.Cells.Item(171, 18).Value = "Hortaliza"
